$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The first column (A) contained the "GENE" count values (4, 14) with the
# header cell A1 empty. Deleting the entire column shifts B:F left to A:E,
# matching the target layout (QS_Astral15 | FNRATE_PHYLONET | TAXON |
# MODEL_CONDITION | GENE) while the old GENE-count column becomes the new
# last column E.
$ws.Columns("A").Delete()
